# PENDIENTES 2025-04-10.xlsx — "nuevos OK resueltos en excel"
#
# Several pending items got marked as resolved by writing "OK" in column C
# (next to the item's description in column B), using the same green/bold
# "OK" look already used on rows 17 and 18 (style: bold 14pt font on a
# solid green fill).  We replicate that by copying the existing C17 cell
# format onto the newly-resolved rows, then filling in the "OK" text, and
# bumping each affected row's height to match (18.75pt, same as rows 17/18).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose pending item just got resolved ("OK" added in column C).
$resolvedRows = @(20, 21, 23, 36, 43, 44, 49, 55)

# Use the existing "OK" cell (C17) as the format template — bold 14pt font
# on green fill — so the new cells reuse the same cell style instead of
# creating a new one.
$ws.Range("C17").Copy()

foreach ($r in $resolvedRows) {
    $cell = $ws.Range("C$r")
    $cell.PasteSpecial(-4122)   # xlPasteFormats — copy formatting only
    $cell.Value = "OK"
    $ws.Rows("${r}:${r}").RowHeight = 18.75
}

$excel.CutCopyMode = $false

# Reflect where the user ended up after making these edits.
$ws.Range("B63").Select() | Out-Null
